$wb = $excel.ActiveWorkbook

# --- Sheet "Recommandations" ---
$ws1 = $wb.Worksheets.Item("Recommandations")

$data1 = @(
    @("NEI-CEDA CI", 0, 5, 4830, 945, "🟡 Observer", "➖ Neutre"),
    @("BRVM - SERVICES PUBLICS", 0, 10, 4212.41, 116.16, "🟡 Observer", "➖ Neutre"),
    @("AIR LIQUIDE CI", 0, 5, 3475, 695, "🟡 Observer", "➖ Neutre"),
    @("BRVM - AUTRES SECTEURS", 0, 5, 2982.72, 597.42, "🟡 Observer", "➖ Neutre"),
    @("BRVM - DISTRIBUTION", 0, 5, 2494.7, 501.88, "🟡 Observer", "➖ Neutre"),
    @("BRVM - TRANSPORT", 0, 5, 1811.32, 361.04, "🟡 Observer", "➖ Neutre"),
    @("BRVM - AGRICULTURE", 0, 5, 1661.37, 332.9, "🟡 Observer", "➖ Neutre"),
    @("BRVM - CONSOMMATION DISCRETIONNAIRE", 0, 5, 861.49, 170.89, "🟡 Observer", "➖ Neutre"),
    @("BRVM - FINANCES", 0, 5, 727.22, 145.87, "🟡 Observer", "➖ Neutre"),
    @("BRVM - SERVICES FINANCIERS", 0, 5, 714.6900000000001, 143.36, "🟡 Observer", "➖ Neutre"),
    @("BRVM-PRESTIGE", 0, 5, 709.14, 140.64, "🟡 Observer", "➖ Neutre"),
    @("BRVM - INDUSTRIELS", 0, 5, 616.54, 122.84, "🟡 Observer", "➖ Neutre"),
    @("BRVM - ENERGIE", 0, 5, 551.52, 113.17, "🟡 Observer", "➖ Neutre"),
    @("BRVM - INDUSTRIE                 (**)", 0, 2, 527.4299999999999, 260.8, "🟡 Observer", "➖ Neutre"),
    @("BRVM - TELECOMMUNICATIONS", 0, 5, 472.03, 95.12, "🟡 Observer", "➖ Neutre"),
    @("BRVM - CONSOMMATION DE BASE          (**)", 0, 2, 442.18, 219.83, "🟡 Observer", "➖ Neutre"),
    @("BRVM-PRINCIPAL                    (**)", 0, 2, 437.2, 218.98, "🟡 Observer", "➖ Neutre"),
    @("BRVM - INDUSTRIE                  (**)", 0, 1, 263.26, 263.26, "🟡 Observer", "➖ Neutre"),
    @("BRVM-PRINCIPAL                   (**)", 0, 1, 218.48, 218.48, "🟡 Observer", "➖ Neutre"),
    @("BRVM - CONSOMMATION DE BASE         (**)", 0, 1, 217.81, 217.81, "🟡 Observer", "➖ Neutre"),
    @("UNILEVER CI (UNLC)", 3, 0, 22.48, 7.49, "🟢 Achat", "✅ Renforcer"),
    @("VIVO ENERGY CI (SHEC)", 3, 0, 19.17, 7.47, "🟢 Achat", "✅ Renforcer"),
    @("SICOR CI (SICC)", 2, 0, 14.82, 7.37, "🟡 Observer", "➖ Neutre"),
    @("SOLIBRA CI (SLBC)", 2, 1, 9.619999999999999, -5.25, "🟡 Observer", "👀 À surveiller"),
    @("ECOBANK TRANS. INCORP. TG (ETIT)", 2, 0, 8.9, 4.35, "🟡 Observer", "➖ Neutre"),
    @("SAFCA CI (SAFC)", 2, 2, 8.41, 7.36, "🟡 Observer", "👀 À surveiller"),
    @("CORIS BANK INTERNATIONAL (CBIBF)", 1, 1, 3.01, -2.7, "🟡 Observer", "👀 À surveiller"),
    @("UNIWAX CI (UNXC)", 2, 1, 1.9, 2.8, "🟡 Observer", "👀 À surveiller"),
    @("NESTLE CI (NTLC)", 2, 1, 0.29, 3.04, "🟡 Observer", "👀 À surveiller"),
    @("SMB CI (SMBC)", 1, 1, 0.21, 4.07, "🟡 Observer", "👀 À surveiller"),
    @("ORAGROUP TOGO (ORGT)", 1, 1, 0.17, 4.29, "🟡 Observer", "👀 À surveiller"),
    @("TOTAL", 0, 4, 0, 0, "🟡 Observer", "➖ Neutre"),
    @("BICI CI (BICC)", 1, 1, -0.33, 3.85, "🟡 Observer", "👀 À surveiller"),
    @("AFRICA GLOBAL LOGISTICS CI (SDSC)", 0, 1, -1.67, -1.67, "🟡 Observer", "➖ Neutre"),
    @("BANK OF AFRICA SENEGAL (BOAS)", 0, 1, -1.7, -1.7, "🟡 Observer", "➖ Neutre"),
    @("CFAO MOTORS CI (CFAC)", 0, 1, -2.31, -2.31, "🟡 Observer", "➖ Neutre"),
    @("SOCIETE GENERALE COTE D'IVOIRE (SGBC)", 0, 1, -2.57, -2.57, "🟡 Observer", "➖ Neutre"),
    @("TRACTAFRIC MOTORS CI (PRSC)", 0, 1, -2.73, -2.73, "🟡 Observer", "➖ Neutre"),
    @("SOCIETE IVOIRIENNE DE BANQUE  (SIBC)", 0, 1, -3.51, -3.51, "🟡 Observer", "➖ Neutre"),
    @("SUCRIVOIRE (SCRC)", 0, 1, -3.85, -3.85, "🟡 Observer", "➖ Neutre"),
    @("NEI-CEDA CI (NEIC)", 1, 2, -4.17, -3.08, "🟡 Observer", "👀 À surveiller"),
    @("ECOBANK COTE D''IVOIRE (ECOC)", 0, 1, -4.38, -4.38, "🟡 Observer", "➖ Neutre"),
    @("SETAO CI (STAC)", 0, 2, -8.52, -4.35, "🟡 Observer", "➖ Neutre"),
    @("SICABLE CI (CABC)", 1, 3, -14.36, 7.37, "🔴 Vente", "⚠️ Risque de décrochage")
)

$r = 2
foreach ($row in $data1) {
    for ($c = 0; $c -lt 7; $c++) {
        $ws1.Cells.Item($r, $c + 1).Value = $row[$c]
    }
    $r++
}

# --- Sheet "Top_YTD" ---
$ws2 = $wb.Worksheets.Item("Top_YTD")

$data2 = @(
    181013923.49,
    13746833.66,
    3175193.6,
    1639400.55,
    770615.64,
    210935.54,
    150833.72,
    14869.52,
    8807.549999999999,
    8362.01
)

$r = 2
foreach ($val in $data2) {
    $ws2.Cells.Item($r, 2).Value = $val
    $r++
}

"Update complete"